$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns J (10) and K (11): widths (closest achievable via character-width rounding) ---
$ws.Columns.Item(10).ColumnWidth = 20.666666666666668
$ws.Columns.Item(11).ColumnWidth = 24

# --- Glossary header cell J1 ---
$ws.Range("J1").Value = "Glossary"
$ws.Range("J1").Borders.Item(7).LineStyle = 1
$ws.Range("J1").HorizontalAlignment = -4108

# --- Left border (thin) down the whole J column box, rows 2-21 ---
$ws.Range("J2:J21").Borders.Item(7).LineStyle = 1

# --- Right border (thin) for K1 (header) with center alignment ---
$ws.Range("K1").Borders.Item(10).LineStyle = 1
$ws.Range("K1").HorizontalAlignment = -4108

# --- Right border (thin) down the whole K column box, rows 2-21 ---
$ws.Range("K2:K21").Borders.Item(10).LineStyle = 1

# --- Glossary entries ---
$ws.Range("J3").Value = "IC"
$ws.Range("K3").Value = "Integrated Circuit"
$ws.Range("J4").Value = "MCU"
$ws.Range("K4").Value = "MicroControler Unit"
$ws.Range("J5").Value = "SBC"
$ws.Range("K5").Value = "Single board computer"

# --- Merge the header cell across J1:K1 (done last so it doesn't disturb the style table order) ---
$ws.Range("J1:K1").Merge() | Out-Null

# --- Expand the table "Tabela2" to cover the extra rows reserved for future entries ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:H450")) | Out-Null

# --- Update the active selection to match the author's last cursor position ---
$ws.Range("J27").Select() | Out-Null
